# Apply edits described by the commit:
# "Removido consulta CcsXi. Documentação. Limpeza de codigo. Ajustes finos finais."
#
# Changes target the "Descrição Tecnica" worksheet:
#  - Header row (row 1) cell order/content updated
#  - Several "Tecnologias Utilizadas" (column D) descriptions updated to mention
#    "Banco remoto de Dados" instead of "Banco de Dados"
#  - Row 12 (C12/D12) values changed to the z-memoria-massa / JCO-only technology text
#  - Row 13 gets a new B13 description, and C13/D13 updated (billing-det / JCO tech text)
#  - Selection moved from D2 to B14

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Descrição Tecnica")
$ws.Activate()

# --- Header row ---
$ws.Range("A1").Value = "Fonte Abap"
$ws.Range("B1").Value = "Descrição da Funcionalidade"
$ws.Range("C1").Value = "Pacote Java gerado"
$ws.Range("D1").Value = "Tecnologias Utilizadas"

# --- Technology description texts ---
$textEjbJaxB   = "Java Dom (W3C) para manipular nodes e Xml Transform para renderizar conteudo. Componente Ejb 3.0 para acesso ao Banco remoto de Dados. Api JaxB para fazer o mapeamento de Xml para Objeto Java."
$textEjbOnly   = "Java Dom (W3C) para manipular nodes e Xml Transform para renderizar conteudo. Componente Ejb 3.0 para acesso ao Banco remoto de Dados."
$textEjbJco    = "Java Dom (W3C) para manipular nodes e Xml Transform para renderizar conteudo. Componente Ejb 3.0 para acesso ao Banco remoto de Dados. JCO para acesso a funções remotas."
$textJcoOnly   = "Java Dom (W3C) para manipular nodes e Xml Transform para renderizar conteudo. Componente JCO para acesso a funções remotas."
$textMemMassa  = "z-memoria-massa"
$textBillingDet = "billing-det"
$textPromove   = "Promove a manipulação padrão de nodes xml, buscando uma data pré definida e atribuindo seu valor aos nodes marcados como pendentes de de atualização."

$ws.Range("D6").Value = $textEjbJaxB
$ws.Range("D7").Value = $textEjbOnly
$ws.Range("D8").Value = $textEjbOnly
$ws.Range("D9").Value = $textEjbJco
$ws.Range("D10").Value = $textEjbJco
$ws.Range("D11").Value = $textEjbJco

$ws.Range("C12").Value = $textMemMassa
$ws.Range("D12").Value = $textJcoOnly

$ws.Range("B13").Value = $textPromove
$ws.Range("C13").Value = $textBillingDet
$ws.Range("D13").Value = $textEjbJco

# --- Update selection to B14 ---
$ws.Range("B14").Select()
